# Fgf1-Fgfr4 LR-pairs sheet update: a new "Resolving-Mac" sending cluster
# was added (cross product against the 4 existing target clusters), and the
# TPM-derived numeric columns were recomputed for every row as a result.
# Row 1 (the header) is unchanged; rows 2-17 are rewritten in full below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @('ECs', 'Fgf1', 'Fgfr4', 'ECs', 3, 1, 1.190640666666667, 3.571922, 0.1136540143525372, 0.1136540143525372, 2, 0.6666666666666666, 0.3659943333333333, 1.097983, 0.006726051721149161, 0.006726051721149162, 0.4357677370362222, 3.921909633326, 0.0007644427788513947, 0.0007644427788513948),
    @('ECs', 'Fgf1', 'Fgfr4', 'FAPs', 3, 1, 1.190640666666667, 3.571922, 0.1136540143525372, 0.1136540143525372, 1, 0.3333333333333333, 0.09131133333333334, 0.273934, 0.001678071748088335, 0.001678071748088335, 0.1087189867942222, 0.9784708811480001, 0.0001907195905418189, 0.0001907195905418189),
    @('ECs', 'Fgf1', 'Fgfr4', 'MuSCs', 3, 1, 1.190640666666667, 3.571922, 0.1136540143525372, 0.1136540143525372, 3, 1, 53.897087, 161.691261, 0.9904923704135933, 0.9904923704135934, 64.17206359707133, 577.5485723736421, 0.1125734340830652, 0.1125734340830652),
    @('ECs', 'Fgf1', 'Fgfr4', 'Resolving-Mac', 3, 1, 1.190640666666667, 3.571922, 0.1136540143525372, 0.1136540143525372, 2, 0.6666666666666666, 0.06004666666666667, 0.18014, 0.001103506117169219, 0.001103506117169219, 0.07149400323111112, 0.6434460290800001, 0.000125417900078863, 0.000125417900078863),
    @('FAPs', 'Fgf1', 'Fgfr4', 'ECs', 3, 1, 5.168173666666667, 15.504521, 0.4933341355895272, 0.4933341355895272, 2, 0.6666666666666666, 0.3659943333333333, 1.097983, 0.006726051721149161, 0.006726051721149162, 1.891522275682556, 17.023700481143, 0.003318190911783573, 0.003318190911783574),
    @('FAPs', 'Fgf1', 'Fgfr4', 'FAPs', 3, 1, 5.168173666666667, 15.504521, 0.4933341355895272, 0.4933341355895272, 1, 0.3333333333333333, 0.09131133333333334, 0.273934, 0.001678071748088335, 0.001678071748088335, 0.4719128284015556, 4.247215455614, 0.0008278500753003658, 0.0008278500753003658),
    @('FAPs', 'Fgf1', 'Fgfr4', 'MuSCs', 3, 1, 5.168173666666667, 15.504521, 0.4933341355895272, 0.4933341355895272, 3, 1, 53.897087, 161.691261, 0.9904923704135933, 0.9904923704135934, 278.5495057434424, 2506.945551690981, 0.4886436973660118, 0.4886436973660119),
    @('FAPs', 'Fgf1', 'Fgfr4', 'Resolving-Mac', 3, 1, 5.168173666666667, 15.504521, 0.4933341355895272, 0.4933341355895272, 2, 0.6666666666666666, 0.06004666666666667, 0.18014, 0.001103506117169219, 0.001103506117169219, 0.3103316014377778, 2.792984412940001, 0.000544397236431432, 0.0005443972364314321),
    @('MuSCs', 'Fgf1', 'Fgfr4', 'ECs', 3, 1, 4.037194, 12.111582, 0.3853751326204581, 0.3853751326204581, 2, 0.6666666666666666, 0.3659943333333333, 1.097983, 0.006726051721149161, 0.006726051721149162, 1.477590126567333, 13.298311139106, 0.002592053074049919, 0.002592053074049919),
    @('MuSCs', 'Fgf1', 'Fgfr4', 'FAPs', 3, 1, 4.037194, 12.111582, 0.3853751326204581, 0.3853751326204581, 1, 0.3333333333333333, 0.09131133333333334, 0.273934, 0.001678071748088335, 0.001678071748088335, 0.3686415670653334, 3.317774103588, 0.0006466871224661861, 0.0006466871224661861),
    @('MuSCs', 'Fgf1', 'Fgfr4', 'MuSCs', 3, 1, 4.037194, 12.111582, 0.3853751326204581, 0.3853751326204581, 3, 1, 53.897087, 161.691261, 0.9904923704135933, 0.9904923704135934, 217.592996253878, 1958.336966284902, 0.3817111286076904, 0.3817111286076905),
    @('MuSCs', 'Fgf1', 'Fgfr4', 'Resolving-Mac', 3, 1, 4.037194, 12.111582, 0.3853751326204581, 0.3853751326204581, 2, 0.6666666666666666, 0.06004666666666667, 0.18014, 0.001103506117169219, 0.001103506117169219, 0.2424200423866667, 2.18178038148, 0.0004252638162515744, 0.0004252638162515744),
    @('Resolving-Mac', 'Fgf1', 'Fgfr4', 'ECs', 1, 0.3333333333333333, 0.08000233333333333, 0.240007, 0.007636717437477471, 0.007636717437477472, 2, 0.6666666666666666, 0.3659943333333333, 1.097983, 0.006726051721149161, 0.006726051721149162, 0.02928040065344444, 0.263523605881, 0.00005136495646427516, 0.00005136495646427517),
    @('Resolving-Mac', 'Fgf1', 'Fgfr4', 'FAPs', 1, 0.3333333333333333, 0.08000233333333333, 0.240007, 0.007636717437477471, 0.007636717437477472, 1, 0.3333333333333333, 0.09131133333333334, 0.273934, 0.001678071748088335, 0.001678071748088335, 0.007305119726444444, 0.065746077538, 0.00001281495977996449, 0.0000128149597799645),
    @('Resolving-Mac', 'Fgf1', 'Fgfr4', 'MuSCs', 1, 0.3333333333333333, 0.08000233333333333, 0.240007, 0.007636717437477471, 0.007636717437477472, 3, 1, 53.897087, 161.691261, 0.9904923704135933, 0.9904923704135934, 4.311892719869666, 38.807034478827, 0.007564110356825882, 0.007564110356825884),
    @('Resolving-Mac', 'Fgf1', 'Fgfr4', 'Resolving-Mac', 1, 0.3333333333333333, 0.08000233333333333, 0.240007, 0.007636717437477471, 0.007636717437477472, 2, 0.6666666666666666, 0.06004666666666667, 0.18014, 0.001103506117169219, 0.001103506117169219, 0.004803873442222223, 0.04323486098, 0.000008427164407349228, 0.000008427164407349231)
)

$nRows = 16
$nCols = 20
$data = New-Object 'object[,]' $nRows,$nCols
for ($r = 0; $r -lt $nRows; $r++) {
  for ($c = 0; $c -lt $nCols; $c++) {
    $data[$r,$c] = $rows[$r][$c]
  }
}

$targetRange = $ws.Range($ws.Cells.Item(2,1), $ws.Cells.Item(1 + $nRows, $nCols))
$targetRange.Value2 = $data
